$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.09950268380201061
$ws.Range("C2").Value = 0.09950268380201061
$ws.Range("D2").Value = 38
$ws.Range("E2").Value = 38
$ws.Range("F2").Value = 192
$ws.Range("G2").Value = 192

$ws.Range("B3").Value = 0.5424796504655186
$ws.Range("C3").Value = 0.5424796504655186
$ws.Range("D3").Value = 25
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 76
$ws.Range("G3").Value = 76

$ws.Range("B4").Value = 0.6509361068228677
$ws.Range("C4").Value = 0.6509361068228677
$ws.Range("D4").Value = 16
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 50
$ws.Range("G4").Value = 50

$ws.Range("B5").Value = 0.004638037914405086
$ws.Range("C5").Value = 0.00408726988385714
$ws.Range("D5").Value = 245
$ws.Range("E5").Value = 237
$ws.Range("F5").Value = 271
$ws.Range("G5").Value = 271

$ws.Range("B6").Value = 0.9391719980261437
$ws.Range("C6").Value = 0.9391719980261437
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 9
$ws.Range("G6").Value = 9

$ws.Range("B7").Value = 0.00005411969235495265
$ws.Range("C7").Value = 0.00003938179468674988
$ws.Range("D7").Value = 851
$ws.Range("E7").Value = 794
$ws.Range("F7").Value = 898
$ws.Range("G7").Value = 898

$ws.Range("B8").Value = 0.08946727459308944
$ws.Range("C8").Value = 0.08588199538181801
$ws.Range("D8").Value = 139
$ws.Range("E8").Value = 135
$ws.Range("F8").Value = 172
$ws.Range("G8").Value = 172

$ws.Range("B9").Value = 0.03507900269087775
$ws.Range("C9").Value = 0.03430511105638601
$ws.Range("D9").Value = 275
$ws.Range("E9").Value = 271
$ws.Range("F9").Value = 306
$ws.Range("G9").Value = 306

$ws.Range("B10").Value = 0.3857449803705244
$ws.Range("C10").Value = 0.381939279255351
$ws.Range("D10").Value = 39
$ws.Range("E10").Value = 37
$ws.Range("F10").Value = 128
$ws.Range("G10").Value = 128

$ws.Range("B11").Value = 0.000000002169669594463701
$ws.Range("C11").Value = 0.000000001192174361104931
$ws.Range("D11").Value = 1173
$ws.Range("E11").Value = 1112
$ws.Range("F11").Value = 1486
$ws.Range("G11").Value = 1484

$ws.Range("B12").Value = 0.03923822944299557
$ws.Range("C12").Value = 0.03923822944299557
$ws.Range("D12").Value = 110
$ws.Range("E12").Value = 110
$ws.Range("F12").Value = 312
$ws.Range("G12").Value = 312

$ws.Range("B13").Value = 0.5035560325870317
$ws.Range("C13").Value = 0.5019982358591089
$ws.Range("D13").Value = 88
$ws.Range("E13").Value = 86
$ws.Range("F13").Value = 92
$ws.Range("G13").Value = 92

$ws.Range("B14").Value = 0.7688953153587761
$ws.Range("C14").Value = 0.7617814239549737
$ws.Range("D14").Value = 30
$ws.Range("E14").Value = 25
$ws.Range("F14").Value = 36
$ws.Range("G14").Value = 36

$ws.Range("B15").Value = 0.7894068617535624
$ws.Range("C15").Value = 0.7868686898597245
$ws.Range("D15").Value = 18
$ws.Range("E15").Value = 17
$ws.Range("F15").Value = 32
$ws.Range("G15").Value = 31

$ws.Range("B16").Value = 0.6266129947137067
$ws.Range("C16").Value = 0.6266129947137067
$ws.Range("D16").Value = 62
$ws.Range("E16").Value = 62
$ws.Range("F16").Value = 65
$ws.Range("G16").Value = 65

$ws.Range("B17").Value = 0.5597822348226835
$ws.Range("C17").Value = 0.5597822348226835
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = 18
$ws.Range("F17").Value = 65
$ws.Range("G17").Value = 65
